$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 06:19:07"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "71%"
$ws.Range("I2").Value = "0.2 mm"
$ws.Range("E3").Value = "2026-02-20 06:19:09"
$ws.Range("I3").Value = "1.2 mm"
$ws.Range("N3").Value = "-6.2 °C 5:58 TU"
$ws.Range("E4").Value = "2026-02-20 06:19:12"
$ws.Range("J4").Value = "1018.7 hPa"
$ws.Range("O4").Value = "8.4 °C"
$ws.Range("E5").Value = "2026-02-20 06:19:14"
$ws.Range("E6").Value = "2026-02-20 06:19:17"
$ws.Range("J6").Value = "1018.7 hPa"
$ws.Range("O6").Value = "5.3 °C"
$ws.Range("E7").Value = "2026-02-20 06:19:19"
$ws.Range("J7").Value = "1018.4 hPa"
$ws.Range("N7").Value = "10.2 °C 5:30 TU"
$ws.Range("E8").Value = "2026-02-20 06:19:22"
$ws.Range("J8").Value = "1019.3 hPa"
$ws.Range("M8").Value = "7.6 °C 5:59 TU"
$ws.Range("E9").Value = "2026-02-20 06:19:24"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "45%"
$ws.Range("O9").Value = "12.5 °C"
$ws.Range("E10").Value = "2026-02-20 06:19:27"
$ws.Range("E11").Value = "2026-02-20 06:19:29"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "35%"
$ws.Range("O11").Value = "7.8 °C"
$ws.Range("E12").Value = "2026-02-20 06:19:31"
$ws.Range("O12").Value = "12.2 °C"
$ws.Range("E13").Value = "2026-02-20 06:19:34"
$ws.Range("J13").Value = "1020.2 hPa"
$ws.Range("N13").Value = "4.1 °C 5:57 TU"
$ws.Range("O13").Value = "5.0 °C"
$ws.Range("E14").Value = "2026-02-20 06:19:36"
$ws.Range("O14").Value = "10.2 °C"
$ws.Range("E15").Value = "2026-02-20 06:19:39"
$ws.Range("N15").Value = "11.6 °C 5:50 TU"
$ws.Range("E16").Value = "2026-02-20 06:19:41"
$ws.Range("L16").Value = "81.7 km/h - 227º 5:56 TU"
$ws.Range("O16").Value = "-5.5 °C"
$ws.Range("E17").Value = "2026-02-20 06:19:44"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "45%"
$ws.Range("K17").Value = "-0.1 MJ/m2"
$ws.Range("E18").Value = "2026-02-20 06:19:47"
$ws.Range("J18").Value = "1019.1 hPa"
$ws.Range("N18").Value = "-0.1 °C 5:51 TU"
$ws.Range("O18").Value = "1.7 °C"
$ws.Range("E19").Value = "2026-02-20 06:19:49"
$ws.Range("O19").Value = "1.6 °C"
$ws.Range("E20").Value = "2026-02-20 06:19:52"
$ws.Range("N20").Value = "-6.1 °C 5:59 TU"
$ws.Range("O20").Value = "-5.3 °C"
$ws.Range("E21").Value = "2026-02-20 06:19:54"
$ws.Range("J21").Value = "1020.3 hPa"
$ws.Range("N21").Value = "3.5 °C 5:50 TU"
$ws.Range("O21").Value = "6.2 °C"
$ws.Range("E22").Value = "2026-02-20 06:19:57"
$ws.Range("E23").Value = "2026-02-20 06:19:59"
$ws.Range("I23").Value = "3.4 mm"
$ws.Range("N23").Value = "-7.1 °C 5:58 TU"
$ws.Range("O23").Value = "-6.5 °C"
$ws.Range("E24").Value = "2026-02-20 06:20:02"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "76%"
$ws.Range("O24").Value = "6.2 °C"
$ws.Range("E25").Value = "2026-02-20 06:20:04"
$ws.Range("I25").Value = "4.4 mm"
$ws.Range("M25").Value = "-3.9 °C 5:53 TU"
$ws.Range("O25").Value = "-4.8 °C"
$ws.Range("E26").Value = "2026-02-20 06:20:07"
$ws.Range("J26").Value = "1018.5 hPa"
$ws.Range("O26").Value = "3.3 °C"
$ws.Range("E27").Value = "2026-02-20 06:20:10"
$ws.Range("K27").Value = "-0.1 MJ/m2"
$ws.Range("O27").Value = "-3.0 °C"
$ws.Range("E28").Value = "2026-02-20 06:20:12"
$ws.Range("J28").Value = "1019.7 hPa"
$ws.Range("N28").Value = "0.6 °C 5:54 TU"
$ws.Range("O28").Value = "2.6 °C"
$ws.Range("E29").Value = "2026-02-20 06:20:15"
$ws.Range("E30").Value = "2026-02-20 06:20:17"
$ws.Range("J30").Value = "1018.2 hPa"
$ws.Range("N30").Value = "5.9 °C 5:30 TU"
$ws.Range("O30").Value = "8.8 °C"
$ws.Range("E31").Value = "2026-02-20 06:20:19"
$ws.Range("J31").Value = "1016.8 hPa"
$ws.Range("N31").Value = "9.5 °C 5:58 TU"
$ws.Range("O31").Value = "10.5 °C"
$ws.Range("E32").Value = "2026-02-20 06:20:22"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "90%"
$ws.Range("N32").Value = "1.4 °C 5:57 TU"
$ws.Range("E33").Value = "2026-02-20 06:20:25"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "46%"
$ws.Range("J33").Value = "1019.2 hPa"
$ws.Range("O33").Value = "4.5 °C"
$ws.Range("E34").Value = "2026-02-20 06:20:27"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "63%"
$ws.Range("M34").Value = "-1.2 °C 5:56 TU"
$ws.Range("E35").Value = "2026-02-20 06:20:30"
$ws.Range("J35").Value = "1024.0 hPa"
$ws.Range("N35").Value = "1.7 °C 5:30 TU"
$ws.Range("E36").Value = "2026-02-20 06:20:32"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "43%"
$ws.Range("J36").Value = "1018.7 hPa"
$ws.Range("O36").Value = "13.5 °C"
$ws.Range("E37").Value = "2026-02-20 06:20:35"
$ws.Range("J37").Value = "1021.4 hPa"
$ws.Range("N37").Value = "-0.3 °C 5:31 TU"
$ws.Range("O37").Value = "1.8 °C"
$ws.Range("E38").Value = "2026-02-20 06:20:37"
$ws.Range("N38").Value = "2.1 °C 5:40 TU"
$ws.Range("O38").Value = "4.1 °C"
$ws.Range("E39").Value = "2026-02-20 06:20:40"
$ws.Range("E40").Value = "2026-02-20 06:20:43"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "48%"
$ws.Range("J40").Value = "1021.0 hPa"
$ws.Range("O40").Value = "7.4 °C"
$ws.Range("E41").Value = "2026-02-20 06:20:45"
$ws.Range("J41").Value = "1020.4 hPa"
$ws.Range("N41").Value = "10.1 °C 5:59 TU"
$ws.Range("O41").Value = "10.8 °C"
$ws.Range("E42").Value = "2026-02-20 06:20:48"
$ws.Range("N42").Value = "3.0 °C 5:52 TU"
$ws.Range("O42").Value = "4.5 °C"
$ws.Range("E43").Value = "2026-02-20 06:20:50"
$ws.Range("N43").Value = "-0.1 °C 5:30 TU"
$ws.Range("O43").Value = "1.2 °C"
$ws.Range("E44").Value = "2026-02-20 06:20:53"
$ws.Range("I44").Value = "4.4 mm"
$ws.Range("E45").Value = "2026-02-20 06:20:56"
$ws.Range("J45").Value = "1027.7 hPa"
$ws.Range("N45").Value = "1.6 °C 5:33 TU"
$ws.Range("E46").Value = "2026-02-20 06:20:58"
$ws.Range("J46").Value = "1023.9 hPa"
$ws.Range("N46").Value = "8.6 °C 5:35 TU"
$ws.Range("O46").Value = "9.6 °C"
